$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.977.50"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "'3.388.36"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'142.31"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "'3.966.35"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "'27.78"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000171"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.368.60"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "'61.094.71"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "'6.11"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").Value = "'13.67"
$ws.Range("E19").Value = "  -5.15%  "
$ws.Range("D20").Value = "'8.96"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").Value = "'382.43"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "'74.83"
$ws.Range("E22").Value = "  +2.88%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("E25").Value = "  -5.26%  "
$ws.Range("D26").Value = "'3.522.01"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'8.02"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.16"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'1.41"
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").Value = "'6.98"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'5.04"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "'3.418.49"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").Value = "'0.0771"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "'26.97"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "'2.452.68"
$ws.Range("E47").Value = "  -5.55%  "
$ws.Range("D48").Value = "'23.04"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'6.72"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").Value = "'2.15"
$ws.Range("E51").Value = "  +7.00%  "
